$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New literal values for columns E..T (Ligand-expressing cells ... Edge total expression derived specificity)
# across rows 2..10, per the updated Natmi LR-pairs output (Hras-Cav1).
$values = @(
    @(3, 1, 9.663483666666666, 28.990451, 0.4172798466714015, 0.4172798466714016, 3, 1, 505.8908573333333, 1517.672572, 0.7018211771568338, 0.7018211771568337, 4888.668036956663, 43998.01233260997, 0.2928558331947461, 0.2928558331947461),
    @(3, 1, 9.663483666666666, 28.990451, 0.4172798466714015, 0.4172798466714016, 3, 1, 88.00803400000001, 264.024102, 0.1220933352041998, 0.1220933352041997, 850.4641990944447, 7654.177791850003, 0.0509470881936085, 0.0509470881936085),
    @(3, 1, 9.663483666666666, 28.990451, 0.4172798466714015, 0.4172798466714016, 3, 1, 126.926974, 380.780922, 0.1760854876389666, 0.1760854876389665, 1226.556740108425, 11039.01066097582, 0.07347692528304693, 0.07347692528304692),
    @(3, 1, 5.297564333333334, 15.892693, 0.2287546509102482, 0.2287546509102482, 3, 1, 505.8908573333333, 1517.672572, 0.7018211771568338, 0.7018211771568337, 2679.989362368488, 24119.9042613164, 0.160544858381931, 0.160544858381931),
    @(3, 1, 5.297564333333334, 15.892693, 0.2287546509102482, 0.2287546509102482, 3, 1, 88.00803400000001, 264.024102, 0.1220933352041998, 0.1220933352041997, 466.2282219651874, 4196.053997686687, 0.02792941827310464, 0.02792941827310463),
    @(3, 1, 5.297564333333334, 15.892693, 0.2287546509102482, 0.2287546509102482, 3, 1, 126.926974, 380.780922, 0.1760854876389666, 0.1760854876389665, 672.4038104003274, 6051.634293602947, 0.04028037425521263, 0.04028037425521262),
    @(3, 1, 8.197232333333334, 24.591697, 0.3539655024183503, 0.3539655024183503, 3, 1, 505.8908573333333, 1517.672572, 0.7018211771568338, 0.7018211771568337, 4146.904892870521, 37322.14403583468, 0.2484204855801567, 0.2484204855801566),
    @(3, 1, 8.197232333333334, 24.591697, 0.3539655024183503, 0.3539655024183503, 3, 1, 88.00803400000001, 264.024102, 0.1220933352041998, 0.1220933352041997, 721.4223018978995, 6492.800717081094, 0.04321682873748662, 0.04321682873748661),
    @(3, 1, 8.197232333333334, 24.591697, 0.3539655024183503, 0.3539655024183503, 3, 1, 126.926974, 380.780922, 0.1760854876389666, 0.1760854876389665, 1040.44989524496, 9364.049057204635, 0.062328188100707, 0.06232818810070698)
)

$cols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $r = $startRow + $i
    $rowVals = $values[$i]
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$r").Value = $rowVals[$j]
    }
}
